$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Total Annual Cost")
$ws1.Range("C4").Value = 14046425.72236278
$ws1.Range("C9").Value = 13885032.54876038
$ws1.Range("C11").Value = 14374259.44341646
$ws1.Range("C12").Value = 13930227.32586067
$ws1.Range("C14").Value = 13583941.30695468
$ws1.Range("C57").Value = 15380436.54137731
$ws1.Range("C71").Value = 14374259.44341646
$ws1.Range("C72").Value = 13885642.04020569
$ws1.Range("C94").Value = 14296162.12402704

$ws2 = $wb.Worksheets.Item("Total Economic Loss")
$ws2.Range("C4").Value = 501.2898235000002
$ws2.Range("C9").Value = 363.3748749999996
$ws2.Range("C11").Value = 781.5028450000001
$ws2.Range("C57").Value = 1641.453089
$ws2.Range("C71").Value = 781.5028450000001
$ws2.Range("C94").Value = 714.7553419999996

$ws10 = $wb.Worksheets.Item("Market Transfers")
$ws10.Range("C4").Value = 1178.593372
$ws10.Range("C9").Value = 1040.674875
$ws10.Range("C11").Value = 1458.802845
$ws10.Range("C12").Value = 121.404176
$ws10.Range("C14").Value = 530.6550559999996
$ws10.Range("C57").Value = 2318.753089
$ws10.Range("C71").Value = 1458.802845
$ws10.Range("C72").Value = 174.0654459999996
$ws10.Range("C94").Value = 1392.055342

$ws11 = $wb.Worksheets.Item("TotalShortage")
$ws11.Range("C4").Value = 677.3035485
$ws11.Range("C9").Value = 677.3000000000001
$ws11.Range("C11").Value = 677.3000000000001
$ws11.Range("C12").Value = 121.404176
$ws11.Range("C14").Value = 530.6550559999996
$ws11.Range("C57").Value = 677.3000000000001
$ws11.Range("C71").Value = 677.3000000000001
$ws11.Range("C72").Value = 174.0654459999996
$ws11.Range("C94").Value = 677.3000000000001
